# End to end test for the 3 fuzzer strategies - bank loan with pin
# Updates the Settings sheet of the Config.xlsx workbook:
#   - Refreshes the local python/user paths to the new machine (adelinas / GithubPhD / TestingTool_v4)
#   - Splits the single "ScriptFullPath" entry into three dedicated entries,
#     one per fuzzer strategy (OfflineAll / DFSSymbolic / Concolic)
#   - Moves "RobotModelFullPath" down below the new entries

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Settings")

# PyhtonHomePath value -> new machine's python install path
$ws.Range("B4").Value = "C:\Users\adelinas\AppData\Local\Programs\Python\Python38"

# ScriptWorkingFolder value -> new repo root
$ws.Range("B6").Value = "C:\GithubPhD\rpa-testing\TestingTool_v4"

# Former "ScriptFullPath" row becomes the "offline-all" strategy entry
$ws.Range("A8").Value = "ScriptFullPathOfflineAll"
$ws.Range("B8").Value = "C:\GithubPhD\rpa-testing\TestingTool_v4\bankLoan_offlineall.py"

# New "DFS symbolic" strategy entry (row 10 used to hold RobotModelFullPath)
$ws.Range("A10").Value = "ScriptFullPathDFSSymbolic"
$ws.Range("B10").Value = "C:\GithubPhD\rpa-testing\TestingTool_v4\bankLoan_dfssymbolic.py"
$ws.Range("C10").Value = ""

# New "Concolic" strategy entry
$ws.Range("A12").Value = "ScriptFullPathConcolic"
$ws.Range("B12").Value = "C:\GithubPhD\rpa-testing\TestingTool_v4\bankLoan_concolic.py"

# RobotModelFullPath moved down to row 14, after the three script-path entries
$ws.Range("A14").Value = "RobotModelFullPath"
$ws.Range("B14").Value = "C:\GithubPhD\rpa-testing\TestingTool_v4\Applications\C#Models\SimpleBankLoanCSharp"
$ws.Range("C14").Value = "*Full path of the robot model under test"

# Leave the cursor on the newly added Concolic path cell, matching the saved view state
$ws.Activate() | Out-Null
$ws.Range("B12").Select() | Out-Null
